$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("High School Units")
$ws.Activate()

# Row 7
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

# Row 8
$ws.Range("A8").Value = 3
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0

# Row 9
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 0

# Row 10
$ws.Range("A10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 1

# Row 11
$ws.Range("A11").Value = 2
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 0

# Row 12
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 1

# Row 13
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 0

# Row 14
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 1

# Row 15
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 0

# Row 16
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 1

# Row 17
$ws.Range("M17").Value = 1
$ws.Range("N17").Value = 0

# Row 18
$ws.Range("A18").Value = 4
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 1

# Row 19
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 0

# Row 20
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 1

# Row 21
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 0

# Row 22
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 1

# Row 23
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 0

# Update the selected cell to A18 to match the saved selection in the file
$ws.Range("A18").Select()
